{"js": "// Apply the \"serah terima dpnnt17 pdh eto/5\" edits to the shoe-label\n// mail-merge table. The document has a single 2-column table; column 1\n// (\"T1\" label) and column 2 (\"T2\" label) each hold the same set of\n// merge-field results that need updating independently.\n\nasync function replaceOnce(scopeBody, searchText, newText) {\n  const results = scopeBody.search(searchText, {\n    matchCase: true,\n    matchWholeWord: true,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst leftCell = table.getCell(0, 0);\nconst rightCell = table.getCell(0, 1);\n\n// Left label (\"T1\" -> \"B1\")\nawait replaceOnce(leftCell.body, \"T1\", \"B1\");\nawait replaceOnce(leftCell.body, \"SUNARTO\", \"NOFARIZAL\");\nawait replaceOnce(leftCell.body, \"42\", \"41\");\nawait replaceOnce(leftCell.body, \"56\", \"57\");\nawait replaceOnce(leftCell.body, \"DP4 NAUTIKA / 33\", \"ETO/5\");\nawait replaceOnce(leftCell.body, \"DP4 NAUTIKA / 33\", \"ETO/5\");\n\n// Right label (\"T2\" -> \"B2\")\nawait replaceOnce(rightCell.body, \"T2\", \"B2\");\nawait replaceOnce(rightCell.body, \"NOOR MAULANA\", \"AGUS LUDI D.S.\");\nawait replaceOnce(rightCell.body, \"41\", \"42\");\nawait replaceOnce(rightCell.body, \"56\", \"57\");\nawait replaceOnce(rightCell.body, \"DP4 NAUTIKA / 33\", \"ETO/5\");\nawait replaceOnce(rightCell.body, \"DP4 NAUTIKA / 33\", \"ETO/5\");\n", "ps1": "# Apply the \"serah terima dpnnt17 pdh eto/5\" edits to the shoe-label\n# mail-merge table. The document has a single 2-column table; column 1\n# (\"T1\" label) and column 2 (\"T2\" label) each hold the same set of\n# merge-field results that need updating independently, so every\n# find/replace is scoped to the owning table cell's Range to avoid\n# cross-contamination between the two (near-identical) labels.\n\n$d = $word.ActiveDocument\n\n# NOTE: wdReplaceAll (2) together with a range re-fetched from a table\n# cell after an earlier edit has shifted the document can mis-locate the\n# match (a quirk of this host's Find implementation), so each call here\n# replaces exactly one occurrence (wdReplaceOne) and loops are unrolled\n# for fields that appear twice per cell (\"DP4 NAUTIKA / 33\").\nfunction Replace-InRange($range, [string]$findText, [string]$replaceText) {\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute(\n        $findText,   # FindText\n        $true,       # MatchCase\n        $true,       # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $replaceText,# ReplaceWith\n        1            # Replace (wdReplaceOne)\n    )\n    if (-not $result) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n$tbl = $d.Tables.Item(1)\n$leftCell = $tbl.Cell(1, 1)\n$rightCell = $tbl.Cell(1, 2)\n\n# Left label (\"T1\" -> \"B1\")\nReplace-InRange $leftCell.Range \"T1\" \"B1\"\nReplace-InRange $leftCell.Range \"SUNARTO\" \"NOFARIZAL\"\nReplace-InRange $leftCell.Range \"42\" \"41\"\nReplace-InRange $leftCell.Range \"56\" \"57\"\nReplace-InRange $leftCell.Range \"DP4 NAUTIKA / 33\" \"ETO/5\"\nReplace-InRange $leftCell.Range \"DP4 NAUTIKA / 33\" \"ETO/5\"\n\n# Right label (\"T2\" -> \"B2\")\nReplace-InRange $rightCell.Range \"T2\" \"B2\"\nReplace-InRange $rightCell.Range \"NOOR MAULANA\" \"AGUS LUDI D.S.\"\nReplace-InRange $rightCell.Range \"41\" \"42\"\nReplace-InRange $rightCell.Range \"56\" \"57\"\nReplace-InRange $rightCell.Range \"DP4 NAUTIKA / 33\" \"ETO/5\"\nReplace-InRange $rightCell.Range \"DP4 NAUTIKA / 33\" \"ETO/5\"\n"}
